$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1914893617021277
$ws.Range("C2").Value = 0.5496453900709219
$ws.Range("J2").Value = 0.007092198581560284
$ws.Range("P2").Value = 0.1382978723404255
$ws.Range("S2").Value = 0.1134751773049645
$ws.Range("C3").Value = 0.03680981595092025
$ws.Range("J3").Value = 0.06134969325153374
$ws.Range("P3").Value = 0.6809815950920245
$ws.Range("S3").Value = 0.2208588957055215
$ws.Range("J4").Value = 0.04
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.26
$ws.Range("B6").Value = 0.06593406593406594
$ws.Range("D6").Value = 0.007326007326007326
$ws.Range("F6").Value = 0.04395604395604396
$ws.Range("J6").Value = 0.271062271062271
$ws.Range("O6").Value = 0.0293040293040293
$ws.Range("Q6").Value = 0.1904761904761905
$ws.Range("R6").Value = 0.0695970695970696
$ws.Range("S6").Value = 0.3223443223443224
$ws.Range("B7").Value = 0.07692307692307693
$ws.Range("D7").Value = 0.01357466063348416
$ws.Range("F7").Value = 0.05429864253393665
$ws.Range("J7").Value = 0.1493212669683258
$ws.Range("O7").Value = 0.01357466063348416
$ws.Range("Q7").Value = 0.1900452488687783
$ws.Range("R7").Value = 0.1040723981900453
$ws.Range("S7").Value = 0.3981900452488688
$ws.Range("B8").Value = 0.07322175732217573
$ws.Range("D8").Value = 0.02510460251046025
$ws.Range("F8").Value = 0.04602510460251046
$ws.Range("J8").Value = 0.1129707112970711
$ws.Range("O8").Value = 0.01882845188284519
$ws.Range("Q8").Value = 0.1861924686192468
$ws.Range("R8").Value = 0.09832635983263599
$ws.Range("S8").Value = 0.4393305439330544
$ws.Range("B9").Value = 0.07224334600760456
$ws.Range("D9").Value = 0.007604562737642586
$ws.Range("F9").Value = 0.07984790874524715
$ws.Range("J9").Value = 0.1064638783269962
$ws.Range("O9").Value = 0.03802281368821293
$ws.Range("Q9").Value = 0.1787072243346008
$ws.Range("R9").Value = 0.07604562737642585
$ws.Range("S9").Value = 0.44106463878327
$ws.Range("B10").Value = 0.1026785714285714
$ws.Range("D10").Value = 0.02455357142857143
$ws.Range("E10").Value = 0.001488095238095238
$ws.Range("F10").Value = 0.08035714285714286
$ws.Range("J10").Value = 0.1101190476190476
$ws.Range("O10").Value = 0.01339285714285714
$ws.Range("Q10").Value = 0.2098214285714286
$ws.Range("R10").Value = 0.07291666666666667
$ws.Range("S10").Value = 0.3846726190476191
$ws.Range("G11").Value = 0.09446254071661238
$ws.Range("J11").Value = 0.1205211726384365
$ws.Range("K11").Value = 0.1335504885993485
$ws.Range("L11").Value = 0.6514657980456026
$ws.Range("G12").Value = 0.7524752475247525
$ws.Range("J12").Value = 0.1930693069306931
$ws.Range("K12").Value = 0.004950495049504951
$ws.Range("L12").Value = 0.01485148514851485
$ws.Range("S12").Value = 0.03465346534653466
$ws.Range("F13").Value = 0.015625
$ws.Range("G13").Value = 0.703125
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.03125
$ws.Range("F15").Value = 0.01224489795918367
$ws.Range("H15").Value = 0.1673469387755102
$ws.Range("I15").Value = 0.08163265306122448
$ws.Range("J15").Value = 0.3387755102040816
$ws.Range("K15").Value = 0.04897959183673469
$ws.Range("M15").Value = 0.01224489795918367
$ws.Range("N15").Value = 0.004081632653061225
$ws.Range("O15").Value = 0.07346938775510205
$ws.Range("S15").Value = 0.2612244897959184
$ws.Range("F16").Value = 0.02777777777777778
$ws.Range("H16").Value = 0.1611111111111111
$ws.Range("I16").Value = 0.09444444444444444
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.1166666666666667
$ws.Range("M16").Value = 0.03888888888888889
$ws.Range("O16").Value = 0.02222222222222222
$ws.Range("S16").Value = 0.1388888888888889
$ws.Range("F17").Value = 0.02946954813359529
$ws.Range("H17").Value = 0.1866404715127701
$ws.Range("I17").Value = 0.1119842829076621
$ws.Range("J17").Value = 0.3516699410609037
$ws.Range("K17").Value = 0.1119842829076621
$ws.Range("M17").Value = 0.02161100196463654
$ws.Range("O17").Value = 0.08055009823182711
$ws.Range("S17").Value = 0.106090373280943
$ws.Range("F18").Value = 0.02926829268292683
$ws.Range("H18").Value = 0.1560975609756098
$ws.Range("I18").Value = 0.1024390243902439
$ws.Range("J18").Value = 0.4146341463414634
$ws.Range("K18").Value = 0.1170731707317073
$ws.Range("M18").Value = 0.02926829268292683
$ws.Range("N18").Value = 0.004878048780487805
$ws.Range("O18").Value = 0.04390243902439024
$ws.Range("S18").Value = 0.1024390243902439
$ws.Range("F19").Value = 0.02062588904694168
$ws.Range("H19").Value = 0.2012802275960171
$ws.Range("I19").Value = 0.1066856330014225
$ws.Range("J19").Value = 0.3563300142247511
$ws.Range("K19").Value = 0.104551920341394
$ws.Range("M19").Value = 0.02702702702702703
$ws.Range("N19").Value = 0.0007112375533428165
$ws.Range("O19").Value = 0.07539118065433854
$ws.Range("S19").Value = 0.1073968705547653
